$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: address -> new text value (values stay as text, matching original inlineStr formatting)
$updates = @{
    "D2" = "298.39"
    "E2" = "0.79%"
    "D3" = "31.32"
    "E3" = "0.66%"
    "E4" = "0.51%"
    "D5" = "0.07944"
    "E5" = "7.84%"
    "D6" = "2.614"
    "E6" = "59.66%"
    "D7" = "7.832"
    "E7" = "1.65%"
    "D8" = "3.831"
    "E8" = "2.22%"
    "D9" = "0.9088"
    "E9" = "-0.88%"
    "D10" = "0.1734"
    "E10" = "3.52%"
    "D11" = "0.07249"
    "E11" = "1.90%"
    "E12" = "0.87%"
    "D13" = "0.03024"
    "E13" = "1.42%"
    "E14" = "0.62%"
    "D15" = "0.001491"
    "E15" = "-0.05%"
    "D16" = "0.006030"
    "E16" = "-2.33%"
    "D17" = "3.504"
    "E17" = "1.63%"
    "E18" = "1.17%"
    "E19" = "0.33%"
    "D20" = "0.1327"
    "E20" = "-0.34%"
    "D21" = "4.632"
    "E21" = "1.72%"
    "E22" = "3.28%"
    "D23" = "0.04582"
    "E23" = "-0.81%"
    "D24" = "0.001260"
    "E24" = "3.67%"
    "D25" = "0.004449"
    "E25" = "0.54%"
    "D26" = "0.0001180"
    "E26" = "-9.00%"
    "D27" = "0.0003431"
    "E27" = "83.28%"
    "D39" = "0.01834"
    "E39" = "8.85%"
    "E40" = "2.72%"
    "D41" = "0.007024"
    "E41" = "-1.76%"
    "D42" = "0.1342"
    "E42" = "1.11%"
    "D43" = "0.002241"
    "E43" = "4.94%"
    "D44" = "0.01043"
    "E44" = "-5.89%"
    "D45" = "0.00006414"
    "E45" = "6.98%"
    "E46" = "0.02%"
    "E47" = "15.29%"
    "D48" = "0.006202"
    "E48" = "-39.28%"
    "D49" = "0.00002101"
    "E49" = "0.02%"
    "D50" = "0.0002001"
    "E50" = "0.09%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}